$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'29.062.99"
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -0.06%  '
$ws.Range('D3').Value = "'1.836.52"
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +0.44%  '
$ws.Range('D4').Value = "'1.0000"
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = "'243.69"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.91%  '
$ws.Range('D6').Value = "'0.6307"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.13%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = "'0.07593"
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +3.66%  '
$ws.Range('D9').Value = "'0.2941"
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +0.38%  '
$ws.Range('D10').Value = "'22.74"
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +0.50%  '
$ws.Range('D11').Value = "'0.07757"
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +1.27%  '
$ws.Range('D12').Value = "'1.827.56"
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -0.07%  '
$ws.Range('D13').Value = "'4.981"
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +0.04%  '
$ws.Range('D14').Value = "'0.6690"
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +1.06%  '
$ws.Range('B15').Value = 'Litecoin'
$ws.Range('C15').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D15').Value = "'83.29"
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +1.45%  '
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').Value = "'0.000009957"
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +14.56%  '
$ws.Range('D17').Value = "'6.108"
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +1.01%  '
$ws.Range('D18').Value = "'29.079.51"
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +0.01%  '
$ws.Range('D19').Value = "'227.21"
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.09%  '
$ws.Range('D20').Value = "'12.47"
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.79%  '
$ws.Range('E21').Value = '  +0.00%  '
$ws.Range('E22').Value = '  +1.50%  '
$ws.Range('D23').Value = "'0.9984"
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -0.26%  '
$ws.Range('D24').Value = "'159.82"
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.99%  '
$ws.Range('D25').Value = "'0.1395"
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +3.09%  '
$ws.Range('D26').Value = "'8.522"
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.95%  '
$ws.Range('D27').Value = "'17.95"
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.48%  '
$ws.Range('E28').Value = '  -0.12%  '
$ws.Range('D29').Value = "'4.114"
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +1.06%  '
$ws.Range('D30').Value = "'4.025"
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +0.29%  '
$ws.Range('D31').Value = "'1.202"
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.17%  '
$ws.Range('D32').Value = "'0.05267"
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -0.60%  '
$ws.Range('D33').Value = "'1.857"
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +1.13%  '
$ws.Range('D34').Value = "'0.7398"
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +1.00%  '
$ws.Range('D35').Value = "'1.139"
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -1.31%  '
$ws.Range('D36').Value = "'2.681"
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +0.94%  '
$ws.Range('D37').Value = "'1.247.53"
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -3.76%  '
$ws.Range('E38').Value = '  +0.81%  '
$ws.Range('D39').Value = "'0.01788"
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +0.01%  '
$ws.Range('D40').Value = "'6.404"
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +1.53%  '
$ws.Range('D41').Value = "'0.9017"
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +0.00%  '
$ws.Range('D42').Value = "'1.002"
$ws.Range('D42').ClearFormats()
$ws.Range('B43').Value = 'Quant'
$ws.Range('C43').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D43').Value = "'102.30"
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.07%  '
$ws.Range('B44').Value = 'RocketPoolETH'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D44').Value = "'1.982.63"
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +0.38%  '
$ws.Range('B45').Value = 'BabyDogeCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D45').Value = "'0.00000000126"
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +4.81%  '
$ws.Range('D46').Value = "'64.59"
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.95%  '
$ws.Range('D47').Value = "'0.5117"
$ws.Range('D47').ClearFormats()
$ws.Range('D48').Value = "'0.4077"
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +2.26%  '
$ws.Range('D49').Value = "'8.982"
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +3.10%  '
$ws.Range('D50').Value = "'0.05767"
$ws.Range('D50').ClearFormats()
$ws.Range('D51').Value = "'6.738"
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +0.79%  '
